$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J5").Value = 0.86
$ws.Range("G6").Value = 6.429
$ws.Range("H8").Value = 5.773
$ws.Range("J8").Value = 1.82
$ws.Range("G9").Value = 7.208
$ws.Range("H9").Value = 6.039
$ws.Range("J9").Value = 0
$ws.Range("G11").Value = 7.487
$ws.Range("H11").Value = 7.513
$ws.Range("J11").Value = 4.39
$ws.Range("G12").Value = 5.3
$ws.Range("H12").Value = 7.5
$ws.Range("G14").Value = 9.664
$ws.Range("H14").Value = 7.047
$ws.Range("J14").Value = 5.23
$ws.Range("J17").Value = 15.46
$ws.Range("G20").Value = 9.353999999999999
$ws.Range("H20").Value = 10.063
$ws.Range("J20").Value = 11.94
$ws.Range("G21").Value = 8.888
$ws.Range("H21").Value = 10.011
